# Update automatic: dades i banners [2026-02-11 19:20]
# Applies cell-value updates to the meteocat daily summary sheet while
# preserving each cell's existing style (General format, border-only xf).
#
# Directly assigning Range.Value with a string that LOOKS like a pure
# number or percentage (e.g. "80%") makes Excel auto-detect it as a
# numeric/percent value and silently switches the cell's number format
# (creating a brand-new style), which the source workbook does not want -
# every touched cell must keep style index 3 (inlineStr, General number
# format). To force literal text without Excel's smart-parsing, each
# value is written as a text formula (="value") and then flattened back
# to a plain value via Copy + PasteSpecial(xlPasteValues, i.e. -4163),
# which keeps General format / the original style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteValues = -4163

function Set-CellText($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
}

Set-CellText $ws.Range("E2") '2026-02-11 19:18:19'
Set-CellText $ws.Range("H2") '80%'
Set-CellText $ws.Range("I2") '6.6 mm'
Set-CellText $ws.Range("O2") '3.1 °C'
Set-CellText $ws.Range("E3") '2026-02-11 19:18:21'
Set-CellText $ws.Range("H3") '81%'
Set-CellText $ws.Range("I3") '2.8 mm'
Set-CellText $ws.Range("E4") '2026-02-11 19:18:24'
Set-CellText $ws.Range("J4") '1002.4 hPa'
Set-CellText $ws.Range("E5") '2026-02-11 19:18:26'
Set-CellText $ws.Range("I5") '2.7 mm'
Set-CellText $ws.Range("O5") '0.6 °C'
Set-CellText $ws.Range("E6") '2026-02-11 19:18:29'
Set-CellText $ws.Range("J6") '1003.0 hPa'
Set-CellText $ws.Range("E7") '2026-02-11 19:18:31'
Set-CellText $ws.Range("H7") '43%'
Set-CellText $ws.Range("I7") '0.1 mm'
Set-CellText $ws.Range("J7") '1003.6 hPa'
Set-CellText $ws.Range("N7") '15.9 °C 18:59 TU'
Set-CellText $ws.Range("O7") '19.2 °C'
Set-CellText $ws.Range("E8") '2026-02-11 19:18:33'
Set-CellText $ws.Range("H8") '54%'
Set-CellText $ws.Range("I8") '0.6 mm'
Set-CellText $ws.Range("N8") '11.5 °C 18:57 TU'
Set-CellText $ws.Range("O8") '15.3 °C'
Set-CellText $ws.Range("E9") '2026-02-11 19:18:36'
Set-CellText $ws.Range("E10") '2026-02-11 19:18:38'
Set-CellText $ws.Range("H10") '75%'
Set-CellText $ws.Range("L10") '24.5 km/h - 224º 18:59 TU'
Set-CellText $ws.Range("O10") '13.6 °C'
Set-CellText $ws.Range("E11") '2026-02-11 19:18:41'
Set-CellText $ws.Range("E12") '2026-02-11 19:18:43'
Set-CellText $ws.Range("E13") '2026-02-11 19:18:45'
Set-CellText $ws.Range("I13") '0.4 mm'
Set-CellText $ws.Range("J13") '1005.2 hPa'
Set-CellText $ws.Range("E14") '2026-02-11 19:18:48'
Set-CellText $ws.Range("N14") '15.2 °C 18:59 TU'
Set-CellText $ws.Range("O14") '19.1 °C'
Set-CellText $ws.Range("E15") '2026-02-11 19:18:50'
Set-CellText $ws.Range("H15") '83%'
Set-CellText $ws.Range("E16") '2026-02-11 19:18:52'
Set-CellText $ws.Range("H16") '65%'
Set-CellText $ws.Range("I16") '6.9 mm'
Set-CellText $ws.Range("E17") '2026-02-11 19:18:55'
Set-CellText $ws.Range("I17") '0.1 mm'
Set-CellText $ws.Range("E18") '2026-02-11 19:18:57'
Set-CellText $ws.Range("H18") '70%'
Set-CellText $ws.Range("J18") '1003.0 hPa'
Set-CellText $ws.Range("O18") '14.2 °C'
Set-CellText $ws.Range("E19") '2026-02-11 19:19:00'
Set-CellText $ws.Range("E20") '2026-02-11 19:19:02'
Set-CellText $ws.Range("I20") '1.2 mm'
Set-CellText $ws.Range("E21") '2026-02-11 19:19:04'
Set-CellText $ws.Range("I21") '2.1 mm'
Set-CellText $ws.Range("E22") '2026-02-11 19:19:07'
Set-CellText $ws.Range("E23") '2026-02-11 19:19:09'
Set-CellText $ws.Range("I23") '4.8 mm'
Set-CellText $ws.Range("E24") '2026-02-11 19:19:11'
Set-CellText $ws.Range("H24") '74%'
Set-CellText $ws.Range("I24") '8.1 mm'
Set-CellText $ws.Range("J24") '1007.0 hPa'
Set-CellText $ws.Range("N24") '11.0 °C 18:35 TU'
Set-CellText $ws.Range("O24") '13.3 °C'
Set-CellText $ws.Range("E25") '2026-02-11 19:19:14'
Set-CellText $ws.Range("H25") '63%'
Set-CellText $ws.Range("I25") '1.8 mm'
Set-CellText $ws.Range("E26") '2026-02-11 19:19:16'
Set-CellText $ws.Range("J26") '1002.9 hPa'
Set-CellText $ws.Range("O26") '7.1 °C'
Set-CellText $ws.Range("E27") '2026-02-11 19:19:19'
Set-CellText $ws.Range("I27") '1.7 mm'
Set-CellText $ws.Range("E28") '2026-02-11 19:19:21'
Set-CellText $ws.Range("H28") '80%'
Set-CellText $ws.Range("J28") '1003.2 hPa'
Set-CellText $ws.Range("O28") '11.0 °C'
Set-CellText $ws.Range("E29") '2026-02-11 19:19:23'
Set-CellText $ws.Range("E30") '2026-02-11 19:19:26'
Set-CellText $ws.Range("J30") '1003.1 hPa'
Set-CellText $ws.Range("E31") '2026-02-11 19:19:28'
Set-CellText $ws.Range("J31") '1002.3 hPa'
Set-CellText $ws.Range("E32") '2026-02-11 19:19:31'
Set-CellText $ws.Range("I32") '3.5 mm'
Set-CellText $ws.Range("E33") '2026-02-11 19:19:33'
Set-CellText $ws.Range("I33") '1.7 mm'
Set-CellText $ws.Range("O33") '6.7 °C'
Set-CellText $ws.Range("E34") '2026-02-11 19:19:36'
Set-CellText $ws.Range("E35") '2026-02-11 19:19:38'
Set-CellText $ws.Range("J35") '1007.6 hPa'
Set-CellText $ws.Range("M35") '14.0 °C 18:59 TU'
Set-CellText $ws.Range("O35") '10.9 °C'
Set-CellText $ws.Range("E36") '2026-02-11 19:19:41'
Set-CellText $ws.Range("J36") '1003.3 hPa'
Set-CellText $ws.Range("E37") '2026-02-11 19:19:43'
Set-CellText $ws.Range("H37") '80%'
Set-CellText $ws.Range("J37") '1004.4 hPa'
Set-CellText $ws.Range("L37") '49.0 km/h - 246º 18:44 TU'
Set-CellText $ws.Range("O37") '9.5 °C'
Set-CellText $ws.Range("E38") '2026-02-11 19:19:46'
Set-CellText $ws.Range("E39") '2026-02-11 19:19:48'
Set-CellText $ws.Range("O39") '1.2 °C'
Set-CellText $ws.Range("E40") '2026-02-11 19:19:50'
Set-CellText $ws.Range("I40") '3.7 mm'
Set-CellText $ws.Range("J40") '1007.0 hPa'
Set-CellText $ws.Range("E41") '2026-02-11 19:19:52'
Set-CellText $ws.Range("H41") '46%'
Set-CellText $ws.Range("I41") '0.1 mm'
Set-CellText $ws.Range("J41") '1004.8 hPa'
Set-CellText $ws.Range("N41") '14.1 °C 18:59 TU'
Set-CellText $ws.Range("O41") '19.0 °C'
Set-CellText $ws.Range("E42") '2026-02-11 19:19:55'
Set-CellText $ws.Range("E43") '2026-02-11 19:19:57'
Set-CellText $ws.Range("I43") '1.1 mm'
Set-CellText $ws.Range("O43") '13.2 °C'
Set-CellText $ws.Range("E44") '2026-02-11 19:19:59'
Set-CellText $ws.Range("I44") '5.9 mm'
Set-CellText $ws.Range("E45") '2026-02-11 19:20:02'
Set-CellText $ws.Range("I45") '3.8 mm'
Set-CellText $ws.Range("J45") '1005.9 hPa'
Set-CellText $ws.Range("E46") '2026-02-11 19:20:04'
Set-CellText $ws.Range("H46") '59%'
Set-CellText $ws.Range("J46") '1007.3 hPa'
Set-CellText $ws.Range("O46") '17.2 °C'

$excel.CutCopyMode = $false

Write-Host "Applied 122 cell updates (meteocat daily summary refresh)."
